$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3030.537
$ws.Range("I137").Value = 1696.762
$ws.Range("J137").Value = 3879.303
$ws.Range("K137").Value = 5090.286
$ws.Range("L137").Value = 11637.909
$ws.Range("M137").Value = -2540.286
$ws.Range("N137").Value = -16737.909
$ws.Range("H138").Value = 1448321.8
$ws.Range("I138").Value = 4125.1
$ws.Range("K138").Value = 12375.3
$ws.Range("M138").Value = -7235.300000000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2110
$ws.Range("I2").Value = 2179.2727
$ws.Range("J2").Value = 1805.2
$ws.Range("K2").Value = 2179.2727
$ws.Range("L2").Value = 1805.2
$ws.Range("M2").Value = -2066.2727
$ws.Range("N2").Value = -2031.2
$ws.Range("H32").Value = 14624.368
$ws.Range("I32").Value = 12815.471
$ws.Range("K32").Value = 12815.471
$ws.Range("M32").Value = -12528.471
$ws.Range("H74").Value = 3724.5417
$ws.Range("I74").Value = 1793.1875
$ws.Range("J74").Value = 7587.25
$ws.Range("K74").Value = 1793.1875
$ws.Range("L74").Value = 7587.25
$ws.Range("M74").Value = -919.1875
$ws.Range("N74").Value = -9335.25
$ws.Range("H77").Value = 3724.5417
$ws.Range("I77").Value = 1793.1875
$ws.Range("J77").Value = 7587.25
$ws.Range("K77").Value = 8965.9375
$ws.Range("L77").Value = 37936.25
$ws.Range("M77").Value = -4597.9375
$ws.Range("N77").Value = -46672.25
$ws.Range("H116").Value = 2110
$ws.Range("I116").Value = 2179.2727
$ws.Range("J116").Value = 1805.2
$ws.Range("K116").Value = 2179.2727
$ws.Range("L116").Value = 1805.2
$ws.Range("M116").Value = 114.7273
$ws.Range("N116").Value = -6393.2
$ws.Range("H140").Value = 39753.066
$ws.Range("J140").Value = 41193.54
$ws.Range("L140").Value = 41193.54
$ws.Range("N140").Value = -51553.54
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2110
$ws.Range("I3").Value = 2179.2727
$ws.Range("J3").Value = 1805.2
$ws.Range("K3").Value = 2179.2727
$ws.Range("L3").Value = 1805.2
$ws.Range("M3").Value = -2065.2727
$ws.Range("N3").Value = -2033.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5571.4
$ws.Range("I31").Value = 3525
$ws.Range("J31").Value = 5798.778
$ws.Range("K31").Value = 3525
$ws.Range("L31").Value = 5798.778
$ws.Range("M31").Value = -3230
$ws.Range("N31").Value = -6388.778
$ws.Range("H34").Value = 5571.4
$ws.Range("I34").Value = 3525
$ws.Range("J34").Value = 5798.778
$ws.Range("K34").Value = 3525
$ws.Range("L34").Value = 5798.778
$ws.Range("M34").Value = -3323
$ws.Range("N34").Value = -6202.778
$ws.Range("H132").Value = 3008.95
$ws.Range("I132").Value = 2740.5
$ws.Range("J132").Value = 3635.3333
$ws.Range("K132").Value = 8221.5
$ws.Range("L132").Value = 10905.9999
$ws.Range("M132").Value = -5691.5
$ws.Range("N132").Value = -15965.9999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2585062
$ws.Range("I107").Value = 10101398
$ws.Range("J107").Value = 1321.2812
$ws.Range("K107").Value = 30304194
$ws.Range("L107").Value = 3963.8436
$ws.Range("M107").Value = -30302274
$ws.Range("N107").Value = -7803.8436
$ws.Range("H113").Value = 749.907
$ws.Range("I113").Value = 785.3333
$ws.Range("J113").Value = 705.1579
$ws.Range("K113").Value = 2355.9999
$ws.Range("L113").Value = 2115.4737
$ws.Range("M113").Value = -185.9998999999998
$ws.Range("N113").Value = -6455.4737
$ws.Range("H126").Value = 1921.909
$ws.Range("I126").Value = 1272.5714
$ws.Range("J126").Value = 3058.25
$ws.Range("K126").Value = 3817.7142
$ws.Range("L126").Value = 9174.75
$ws.Range("M126").Value = 1122.2858
$ws.Range("N126").Value = -19054.75
$ws.Range("H129").Value = 2153.5715
$ws.Range("I129").Value = 3220
$ws.Range("J129").Value = 1497.3077
$ws.Range("K129").Value = 9660
$ws.Range("L129").Value = 4491.9231
$ws.Range("M129").Value = -4660
$ws.Range("N129").Value = -14491.9231
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6181.5454
$ws.Range("I80").Value = 12930
$ws.Range("K80").Value = 12930
$ws.Range("M80").Value = -11932
$ws.Range("H83").Value = 6181.5454
$ws.Range("I83").Value = 12930
$ws.Range("K83").Value = 64650
$ws.Range("M83").Value = -59658
$ws.Range("H136").Value = 65326
$ws.Range("J136").Value = 65326
$ws.Range("L136").Value = 195978
$ws.Range("N136").Value = -201078
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 6970
$ws.Range("I9").Value = 455
$ws.Range("J9").Value = 20000
$ws.Range("K9").Value = 455
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = -231
$ws.Range("N9").Value = -20448
$ws.Range("H16").Value = 595.8
$ws.Range("I16").Value = 595.8
$ws.Range("K16").Value = 595.8
$ws.Range("M16").Value = -425.8
$ws.Range("H20").Value = 11000
$ws.Range("J20").Value = 11000
$ws.Range("L20").Value = 11000
$ws.Range("N20").Value = -11452
$ws.Range("H22").Value = 2046
$ws.Range("I22").Value = 1952
$ws.Range("J22").Value = 2163.5
$ws.Range("K22").Value = 1952
$ws.Range("L22").Value = 2163.5
$ws.Range("M22").Value = -1657
$ws.Range("N22").Value = -2753.5
$ws.Range("H27").Value = 2046
$ws.Range("I27").Value = 1952
$ws.Range("J27").Value = 2163.5
$ws.Range("K27").Value = 1952
$ws.Range("L27").Value = 2163.5
$ws.Range("M27").Value = -1845
$ws.Range("N27").Value = -2377.5
$ws.Range("H30").Value = 10000
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 10000
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 10000
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -10216
$ws.Range("H46").Value = 650
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 680
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 680
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -1056
$ws.Range("H56").Value = 21833.334
$ws.Range("I56").Value = 11500
$ws.Range("J56").Value = 42500
$ws.Range("K56").Value = 11500
$ws.Range("L56").Value = 42500
$ws.Range("M56").Value = -10809
$ws.Range("N56").Value = -43882
$ws.Range("H82").Value = 1964.8422
$ws.Range("I82").Value = 1170.1538
$ws.Range("J82").Value = 3686.6667
$ws.Range("K82").Value = 1170.1538
$ws.Range("L82").Value = 3686.6667
$ws.Range("M82").Value = -809.1538
$ws.Range("N82").Value = -4408.6667
$ws.Range("H85").Value = 1964.8422
$ws.Range("I85").Value = 1170.1538
$ws.Range("J85").Value = 3686.6667
$ws.Range("K85").Value = 1170.1538
$ws.Range("L85").Value = 3686.6667
$ws.Range("M85").Value = 77.84619999999995
$ws.Range("N85").Value = -6182.6667
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 9166.666999999999
$ws.Range("J26").Value = 11750
$ws.Range("L26").Value = 11750
$ws.Range("N26").Value = -12336
$ws.Range("H29").Value = 5142.857
$ws.Range("I29").Value = 2000
$ws.Range("K29").Value = 2000
$ws.Range("M29").Value = -1710
$ws.Range("H39").Value = 20000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 20000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 20000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -20826
$ws.Range("H81").Value = 3896.25
$ws.Range("I81").Value = 3250
$ws.Range("J81").Value = 4111.6665
$ws.Range("K81").Value = 6500
$ws.Range("L81").Value = 8223.333000000001
$ws.Range("M81").Value = -5439
$ws.Range("N81").Value = -10345.333
$ws.Range("H84").Value = 3896.25
$ws.Range("I84").Value = 3250
$ws.Range("J84").Value = 4111.6665
$ws.Range("K84").Value = 32500
$ws.Range("L84").Value = 41116.665
$ws.Range("M84").Value = -27196
$ws.Range("N84").Value = -51724.665
$ws.Range("H132").Value = 1755.85
$ws.Range("J132").Value = 2528.0625
$ws.Range("L132").Value = 7584.1875
$ws.Range("N132").Value = -12644.1875
$ws.Range("H136").Value = 4427.65
$ws.Range("I136").Value = 2293.7917
$ws.Range("K136").Value = 6881.375100000001
$ws.Range("M136").Value = -4331.375100000001
